$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 - copy formatting (bold/border/alignment) from an existing header cell
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Updated numeric predictions
$ws.Range("B2").Value = 0.4957478377455467
$ws.Range("C2").Value = 0.9901286077698384
$ws.Range("D2").Value = 0.572291725878098

# New model description cell F2 (contains an embedded line break)
$modelo = "Pipeline(steps=[('model'," + [char]10 + "                 AdaBoostRegressor(learning_rate=0.1, n_estimators=150))])"
$ws.Range("F2").Value = $modelo
